$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns I-R: reset the "Итог: <amount>" totals to zero
$ws.Range("I2:R2").Value = "Итог: 0"

# Columns S-T held plain numeric-looking text totals; reset to "0" while
# keeping them as text (leading apostrophe forces text, not a number)
$ws.Range("S2:T2").Value = "'0"
